# Move health facility summary information to a different screen:
# the "distance_to_supply" / "Immunization Services Offered" question block
# used to sit after an extra begin/end-screen pair; that empty screen
# (rows 23:24) is removed so the block moves up onto its own screen.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows("23:24").Delete()

# Reflect the resulting selection/scroll position from the edit.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B22").Select()
